# Refresh coin table cells (prices, 1h volume %, and a few reordered coin rows)
# per the Jan 1 2023 GitHub Actions data-refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'244.20"
$ws.Range("E2").Value = "'-0.74%"
$ws.Range("D3").Value = "'27.20"
$ws.Range("E3").Value = "'4.25%"
$ws.Range("D4").Value = "'5.157"
$ws.Range("E4").Value = "'1.05%"
$ws.Range("E5").Value = "'0.68%"
$ws.Range("D6").Value = "'6.475"
$ws.Range("E6").Value = "'-0.04%"
$ws.Range("D7").Value = "'0.8158"
$ws.Range("E7").Value = "'0.47%"
$ws.Range("D8").Value = "'0.8306"
$ws.Range("E8").Value = "'-1.77%"
$ws.Range("D9").Value = "'0.1329"
$ws.Range("E9").Value = "'-1.17%"
$ws.Range("D10").Value = "'0.06925"
$ws.Range("E10").Value = "'-1.15%"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.03108"
$ws.Range("E11").Value = "'-3.64%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.02932"
$ws.Range("E12").Value = "'5.86%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09389"
$ws.Range("E13").Value = "'-0.12%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001518"
$ws.Range("E14").Value = "'0.16%"
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D15").Value = "'0.04236"
$ws.Range("E15").Value = "'-9.75%"
$ws.Range("B16").Value = "One"
$ws.Range("C16").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D16").Value = "'0.0005957"
$ws.Range("E16").Value = "'-93.93%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.006186"
$ws.Range("E17").Value = "'0.46%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.535"
$ws.Range("E18").Value = "'-0.65%"
$ws.Range("B19").Value = "GateToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D19").Value = "'3.005"
$ws.Range("E19").Value = "'-0.61%"
$ws.Range("B20").Value = "BTSEToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D20").Value = "'2.309"
$ws.Range("E20").Value = "'9.03%"
$ws.Range("B21").Value = "BitpandaEcosystemToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D21").Value = "'0.3113"
$ws.Range("E21").Value = "'-2.20%"
$ws.Range("E22").Value = "'-2.15%"
$ws.Range("D23").Value = "'3.737"
$ws.Range("E23").Value = "'-0.66%"
$ws.Range("E24").Value = "'-0.11%"
$ws.Range("D25").Value = "'0.001225"
$ws.Range("E25").Value = "'-1.89%"
$ws.Range("D26").Value = "'0.004477"
$ws.Range("E26").Value = "'-2.94%"
$ws.Range("D27").Value = "'0.00009797"
$ws.Range("E27").Value = "'2.06%"
$ws.Range("E28").Value = "'-0.49%"
$ws.Range("D40").Value = "'0.03655"
$ws.Range("D41").Value = "'0.1052"
$ws.Range("E41").Value = "'-0.07%"
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D42").Value = "'0.006048"
$ws.Range("E42").Value = "'-1.17%"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002299"
$ws.Range("E43").Value = "'-8.03%"
$ws.Range("D44").Value = "'0.008227"
$ws.Range("E44").Value = "'-5.57%"
$ws.Range("D45").Value = "'0.00005407"
$ws.Range("E45").Value = "'2.11%"
$ws.Range("E46").Value = "'-0.03%"
$ws.Range("D47").Value = "'0.1089"
$ws.Range("E47").Value = "'-17.96%"
$ws.Range("D48").Value = "'0.004715"
$ws.Range("E48").Value = "'130.08%"
$ws.Range("E49").Value = "'-0.03%"
$ws.Range("E50").Value = "'-0.03%"
